$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 166666670
$ws.Range("I46").Value = 166666670
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 500000010
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -499999891
$ws.Range("N46").ClearContents()

$ws.Range("H60").Value = 166666670
$ws.Range("I60").Value = 166666670
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 500000010
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -499999526
$ws.Range("N60").ClearContents()

$ws.Range("H129").Value = 1505.9697
$ws.Range("I129").Value = 526.4545000000001
$ws.Range("J129").Value = 1995.7273
$ws.Range("K129").Value = 1579.3635
$ws.Range("L129").Value = 5987.1819
$ws.Range("M129").Value = 3420.6365
$ws.Range("N129").Value = -15987.1819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5543.2446
$ws.Range("I32").Value = 4716.489
$ws.Range("J32").Value = 17669
$ws.Range("K32").Value = 4716.489
$ws.Range("L32").Value = 17669
$ws.Range("M32").Value = -4429.489
$ws.Range("N32").Value = -18243

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 508
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 16
$ws.Range("K25").Value = 1000
$ws.Range("L25").Value = 16
$ws.Range("M25").Value = -765
$ws.Range("N25").Value = -486

$ws.Range("H54").Value = 3919.75
$ws.Range("I54").Value = 1839.5
$ws.Range("J54").Value = 6000
$ws.Range("K54").Value = 1839.5
$ws.Range("L54").Value = 6000
$ws.Range("M54").Value = -1355.5
$ws.Range("N54").Value = -6968

$ws.Range("H57").Value = 59666.668
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 59666.668
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 59666.668
$ws.Range("N57").Value = -61106.668

$ws.Range("H58").Value = 20923.334
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 20923.334
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 20923.334
$ws.Range("N58").Value = -21511.334

$ws.Range("H60").Value = 20780
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 20780
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 20780
$ws.Range("N60").Value = -21978

$ws.Range("H136").Value = 59666.668
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 59666.668
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 59666.668
$ws.Range("N136").Value = -69866.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2825.2222
$ws.Range("I16").Value = 2763.6
$ws.Range("J16").Value = 3133.3333
$ws.Range("K16").Value = 2763.6
$ws.Range("L16").Value = 3133.3333
$ws.Range("M16").Value = -2476.6
$ws.Range("N16").Value = -3707.3333

$ws.Range("H20").Value = 39499.832
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 39499.832
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 39499.832
$ws.Range("N20").Value = -39971.832

$ws.Range("H30").Value = 39499.832
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 39499.832
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 39499.832
$ws.Range("N30").Value = -39681.832

$ws.Range("H58").Value = 5593.074
$ws.Range("I58").Value = 2753.2727
$ws.Range("J58").Value = 7545.4375
$ws.Range("K58").Value = 2753.2727
$ws.Range("L58").Value = 7545.4375
$ws.Range("M58").Value = -2550.2727
$ws.Range("N58").Value = -7951.4375

$ws.Range("H113").Value = 2825.2222
$ws.Range("I113").Value = 2763.6
$ws.Range("J113").Value = 3133.3333
$ws.Range("K113").Value = 2763.6
$ws.Range("L113").Value = 3133.3333
$ws.Range("M113").Value = -593.5999999999999
$ws.Range("N113").Value = -7473.3333

$ws.Range("H116").Value = 34000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 34000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 34000
$ws.Range("N116").Value = -43178

$ws.Range("H128").Value = 39499.832
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 39499.832
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 39499.832
$ws.Range("N128").Value = -49459.832

$ws.Range("H136").Value = 5593.074
$ws.Range("I136").Value = 2753.2727
$ws.Range("J136").Value = 7545.4375
$ws.Range("K136").Value = 8259.8181
$ws.Range("L136").Value = 22636.3125
$ws.Range("M136").Value = -5709.8181
$ws.Range("N136").Value = -27736.3125

$ws.Range("H140").Value = 47194.75
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 47194.75
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 47194.75
$ws.Range("N140").Value = -57554.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 616.6667
$ws.Range("I28").Value = 550
$ws.Range("J28").Value = 750
$ws.Range("K28").Value = 1650
$ws.Range("L28").Value = 2250
$ws.Range("M28").Value = -1418
$ws.Range("N28").Value = -2714

$ws.Range("H33").Value = 235.8
$ws.Range("I33").Value = 97
$ws.Range("J33").Value = 328.33334
$ws.Range("K33").Value = 582
$ws.Range("L33").Value = 1970.00004
$ws.Range("M33").Value = -299
$ws.Range("N33").Value = -2536.00004

$ws.Range("H34").Value = 901.5
$ws.Range("I34").Value = 99.5
$ws.Range("J34").Value = 1168.8334
$ws.Range("K34").Value = 298.5
$ws.Range("L34").Value = 3506.5002
$ws.Range("M34").Value = -214.5
$ws.Range("N34").Value = -3674.5002

$ws.Range("H45").Value = 1337.375
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1337.375
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 4012.125
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -5076.125

$ws.Range("H98").Value = 1566.2778
$ws.Range("I98").Value = 800
$ws.Range("J98").Value = 2179.3
$ws.Range("K98").Value = 2400
$ws.Range("L98").Value = 6537.900000000001
$ws.Range("M98").Value = -902
$ws.Range("N98").Value = -9533.900000000001

$ws.Range("H106").Value = 4070.4546
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 4070.4546
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 12211.3638
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -14103.3638

$ws.Range("H108").Value = 2778.4443
$ws.Range("I108").Value = 591.8182
$ws.Range("J108").Value = 6214.5713
$ws.Range("K108").Value = 1775.4546
$ws.Range("L108").Value = 18643.7139
$ws.Range("M108").Value = 1104.5454
$ws.Range("N108").Value = -24403.7139

$ws.Range("H115").Value = 2823.5
$ws.Range("I115").Value = 2225.2222
$ws.Range("J115").Value = 3106.8948
$ws.Range("K115").Value = 6675.6666
$ws.Range("L115").Value = 9320.6844
$ws.Range("M115").Value = -5500.6666
$ws.Range("N115").Value = -11670.6844

$ws.Range("H119").Value = 2376.3333
$ws.Range("I119").Value = 2376.3333
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 7128.999899999999
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = -2290.999899999999

$ws.Range("H120").Value = 15009.333
$ws.Range("I120").Value = 5030
$ws.Range("J120").Value = 19999
$ws.Range("K120").Value = 15090
$ws.Range("L120").Value = 59997
$ws.Range("M120").Value = -10252
$ws.Range("N120").Value = -69673

$ws.Range("H129").Value = 2034.6
$ws.Range("I129").Value = 972.75
$ws.Range("J129").Value = 2534.2942
$ws.Range("K129").Value = 2918.25
$ws.Range("L129").Value = 7602.882599999999
$ws.Range("M129").Value = 2081.75
$ws.Range("N129").Value = -17602.8826

$ws.Range("H131").Value = 1315.5333
$ws.Range("I131").Value = 1641.625
$ws.Range("J131").Value = 1196.9546
$ws.Range("K131").Value = 4924.875
$ws.Range("L131").Value = 3590.8638
$ws.Range("M131").Value = 115.125
$ws.Range("N131").Value = -13670.8638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 12000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 12000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 12000
$ws.Range("N15").Value = -12576

$ws.Range("H81").Value = 12000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 12000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 12000
$ws.Range("N81").Value = -13996

$ws.Range("H84").Value = 12000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 12000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 36000
$ws.Range("N84").Value = -45984

$ws.Range("H107").Value = 799.96155
$ws.Range("I107").Value = 814.4666999999999
$ws.Range("J107").Value = 780.1818
$ws.Range("K107").Value = 814.4666999999999
$ws.Range("L107").Value = 780.1818
$ws.Range("M107").Value = 1105.5333
$ws.Range("N107").Value = -4620.1818

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 33641.89
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 33641.89
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 33641.89
$ws.Range("N127").Value = -43561.89

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 6868.8887
$ws.Range("I51").Value = 4606.6665
$ws.Range("J51").Value = 8000
$ws.Range("K51").Value = 4606.6665
$ws.Range("L51").Value = 8000
$ws.Range("M51").Value = -4096.6665
$ws.Range("N51").Value = -9020

$ws.Range("H126").Value = 1619.0488
$ws.Range("I126").Value = 1769.2333
$ws.Range("J126").Value = 1209.4546
$ws.Range("K126").Value = 5307.699900000001
$ws.Range("L126").Value = 3628.3638
$ws.Range("M126").Value = -2837.699900000001
$ws.Range("N126").Value = -8568.363799999999
